# Navi bar fix for smartphone view: expose Mimiku's Instagram link in G2
# (new "twitter" column slot reused for the Instagram hyperlink) and scroll
# the sheet so column G is in view, matching the author's selection change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the Instagram hyperlink to G2. Excel auto-creates the shared string,
# the relationship (rId2) and applies the built-in "Hyperlink" look (underline
# + themed hyperlink colour) to the cell - same as it did for E2's mailto link.
$ws.Hyperlinks.Add($ws.Range("G2"), "https://www.instagram.com/mimiku1210/")

# Move the selection to G2 and scroll column B into the left edge of the
# viewport, mirroring the author's resulting view state.
[void]$ws.Range("B1").Select()
$excel.ActiveWindow.ScrollColumn = 2
[void]$ws.Range("G2").Select()
